$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Murali Vijay / Chennai Super Kings) innings stats
$ws.Range("C2").Value = "10"
$ws.Range("D2").Value = "15"
$ws.Range("E2").Value = "1"

# Row 4 (Murali Vijay / Chennai Super Kings) innings stats
$ws.Range("C4").Value = "21"
$ws.Range("D4").Value = "21"
$ws.Range("E4").Value = "3"
